$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Reorder two pairs of country rows (their labels swap position in the
# shared-string table, while the per-row statistics keep following the
# physical row). Swap the country names so that the correct stats line up
# with the correct country for rows 39/40 and 107/108. ---
$tmp = $ws.Cells.Item(39, 1).Value2
$ws.Cells.Item(39, 1).Value2 = $ws.Cells.Item(40, 1).Value2
$ws.Cells.Item(40, 1).Value2 = $tmp

$tmp = $ws.Cells.Item(107, 1).Value2
$ws.Cells.Item(107, 1).Value2 = $ws.Cells.Item(108, 1).Value2
$ws.Cells.Item(108, 1).Value2 = $tmp

# --- Update the "last updated" timestamp string ---
$ws.Cells.Item(1, 1).Value2 = "Datos actualizados a 28 de Julio de 2020 a las 19:08"

# --- Refresh the per-country statistics (columns B..H) with the newly
# reported numbers. Each tuple is (row, column, newValue). ---
$updates = @(
    @(4, 2, 4455061),
    @(4, 3, 21651),
    @(4, 4, 2141111),
    @(4, 5, 2162402),
    @(4, 7, 473),
    @(4, 8, 151548),
    @(5, 2, 2455905),
    @(5, 3, 12425),
    @(5, 5, 700221),
    @(5, 7, 338),
    @(5, 8, 88017),
    @(6, 2, 1530364),
    @(6, 3, 47861),
    @(6, 4, 987357),
    @(6, 5, 508783),
    @(6, 7, 776),
    @(6, 8, 34224),
    @(11, 2, 349800),
    @(11, 3, 1877),
    @(11, 4, 322332),
    @(11, 5, 18228),
    @(11, 7, 53),
    @(11, 8, 9240),
    @(12, 2, 327690),
    @(12, 3, 1828),
    @(12, 7, 2),
    @(12, 8, 28436),
    @(18, 2, 246488),
    @(18, 3, 181),
    @(18, 5, 12609),
    @(18, 7, 11),
    @(18, 8, 35123),
    @(21, 2, 207508),
    @(21, 3, 129),
    @(21, 5, 6902),
    @(39, 2, 65791),
    @(39, 3, 1806),
    @(39, 4, 32157),
    @(39, 5, 33148),
    @(39, 7, 12),
    @(39, 8, 486),
    @(40, 2, 65149),
    @(40, 3, 770),
    @(40, 4, 55681),
    @(40, 5, 9026),
    @(40, 7, 4),
    @(40, 8, 442),
    @(41, 2, 64690),
    @(41, 3, 534),
    @(41, 4, 32014),
    @(41, 5, 31575),
    @(41, 7, 18),
    @(41, 8, 1101),
    @(60, 2, 28615),
    @(60, 3, 642),
    @(60, 5, 8604),
    @(60, 7, 11),
    @(60, 8, 1174),
    @(61, 2, 25929),
    @(61, 3, 37),
    @(61, 5, 801),
    @(68, 2, 18581),
    @(68, 3, 606),
    @(68, 4, 7908),
    @(68, 5, 10374),
    @(68, 7, 14),
    @(68, 8, 299),
    @(104, 2, 4279),
    @(104, 3, 52),
    @(104, 5, 2702),
    @(104, 7, 1),
    @(104, 8, 203),
    @(105, 2, 4023),
    @(105, 3, 141),
    @(105, 4, 1710),
    @(105, 5, 2259),
    @(105, 7, 3),
    @(105, 8, 54),
    @(107, 2, 3506),
    @(107, 3, 137),
    @(107, 4, 2547),
    @(107, 5, 944),
    @(107, 8, 15),
    @(108, 2, 3439),
    @(108, 4, 2492),
    @(108, 5, 839),
    @(108, 8, 108),
    @(132, 2, 1786),
    @(132, 3, 3),
    @(132, 4, 1336),
    @(132, 5, 384),
    @(134, 2, 1720),
    @(134, 3, 19),
    @(134, 4, 602),
    @(134, 5, 1107),
    @(138, 2, 1468),
    @(138, 3, 13),
    @(138, 4, 1168),
    @(138, 5, 250),
    @(142, 2, 1177),
    @(142, 3, 10),
    @(142, 5, 459),
    @(149, 2, 926),
    @(149, 3, 4),
    @(149, 5, 41),
    @(181, 2, 153),
    @(181, 3, 5),
    @(181, 5, 17)
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value2 = $u[2]
}
